$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")
$wb.ActiveSheet.Activate()

# Update the measured AMS force values in column C (rows 2-9).
# Column D holds =Cn/C5 formulas (and the "AMS %" chart series built on
# top of it) and will recalculate automatically.
$ws.Range("C2").Value = 99.767129999999995
$ws.Range("C3").Value = 267.94510000000002
$ws.Range("C4").Value = 457.38409999999999
$ws.Range("C5").Value = 532.86289999999997
$ws.Range("C6").Value = 1125.6400000000001
$ws.Range("C7").Value = 1097.0319999999999
$ws.Range("C8").Value = 2338.2159999999999
$ws.Range("C9").Value = 1865.671

# Update selected cell in the sheet view.
$ws.Range("C10").Select()

$excel.CalculateFullRebuild()
